$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# CIERRE 27 SEP 23 - record the new entry in row 37 (4 botellones, -212)
$ws.Range("B37").Value = 45195
$ws.Range("C37").Value = "4 botellones"
$ws.Range("D37").Value = -212

# Move the active selection to C38, matching the user's last click location
$ws.Range("C38").Select()
